$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header row: swap E1/F1 so that E1="fecha_registro", F1="costo" ---
$ws.Range("E1").Value = "fecha_registro"
$ws.Range("F1").Value = "costo"

# --- Update row 2 (Aguardiente Amarillo Botella) ---
$ws.Range("C2").Value = 997
$ws.Range("E2").Value = "31/1/2026"
$ws.Range("F2").Value = 45000

# --- Add new row 3 (Cerveza Corona) ---
$ws.Range("A3").Value = "Cerveza Corona"
$ws.Range("B3").Value = "Cervezas"
$ws.Range("C3").Value = 996
$ws.Range("D3").Value = 10000
$ws.Range("E3").Value = "31/1/2026"
$ws.Range("F3").Value = 5000
